$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1: LP1912
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:37:20"
$ws1.Range("A3").Value = "Total filas: 14"

# Update Hora_Scrap column (A) for existing rows 6..15 to new scrape time
for ($r = 6; $r -le 15; $r++) {
    $ws1.Cells.Item($r, 1).Value = "04:37:20"
}

# Row 6: 215_EL PELIGRO
$ws1.Range("B6").Value = "04:47"
$ws1.Range("D6").Value = 10

# Row 7: 11_ETCHEVERRY
$ws1.Range("D7").Value = 16

# Row 8: 17_ROMERO
$ws1.Range("D8").Value = 34

# Row 9: 23_HERNANDEZ
$ws1.Range("B9").Value = "05:22"
$ws1.Range("D9").Value = 45

# Row 10: 81_EL PELIGRO
$ws1.Range("B10").Value = "05:32"
$ws1.Range("D10").Value = 55

# Row 11: 14_ABASTO
$ws1.Range("B11").Value = "05:44"
$ws1.Range("D11").Value = 67

# Row 12: 17_ROMERO
$ws1.Range("B12").Value = "05:52"
$ws1.Range("D12").Value = 75

# Row 13: 16_SANTA ANA
$ws1.Range("B13").Value = "06:01"
$ws1.Range("D13").Value = 84

# Row 14: 10_OLMOS
$ws1.Range("B14").Value = "06:04"
$ws1.Range("D14").Value = 87

# Row 15: 215A_EL PATO
$ws1.Range("B15").Value = "06:11"
$ws1.Range("D15").Value = 94

# New rows 16-19
$ws1.Range("A16").Value = "04:37:20"
$ws1.Range("B16").Value = "06:24"
$ws1.Range("C16").Value = "11_ETCHEVERRY"
$ws1.Range("D16").Value = 107
$ws1.Range("E16").Value = "LP1912"

$ws1.Range("A17").Value = "04:37:20"
$ws1.Range("B17").Value = "06:27"
$ws1.Range("C17").Value = "23_HERNANDEZ"
$ws1.Range("D17").Value = 110
$ws1.Range("E17").Value = "LP1912"

$ws1.Range("A18").Value = "04:37:20"
$ws1.Range("B18").Value = "06:31"
$ws1.Range("C18").Value = "17X38_ROMERO"
$ws1.Range("D18").Value = 114
$ws1.Range("E18").Value = "LP1912"

$ws1.Range("A19").Value = "04:37:20"
$ws1.Range("B19").Value = "06:31"
$ws1.Range("C19").Value = "16_SANTA ANA"
$ws1.Range("D19").Value = 114
$ws1.Range("E19").Value = "LP1912"

# ----------------------------------------------------------------------
# Sheet 2: LP1912-215
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:37:20"

$ws2.Range("A6").Value = "04:37:20"
$ws2.Range("B6").Value = "04:47"
$ws2.Range("D6").Value = 10

$ws2.Range("A7").Value = "04:37:20"
$ws2.Range("B7").Value = "06:11"
$ws2.Range("D7").Value = 94

# ----------------------------------------------------------------------
# Sheet 3: 6203-6173
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:37:20"
